$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","P","Q","R","S","T","U","V","W","X","Z")
foreach ($col in $cols) {
    $ws.Range($col + "4").Value = "np.nan"
}

$ws.Range("A4").Value = "Missing data"

$ws.Range("N17").Select()
